$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.742.92'
$ws.Range('E2').Value = '  +2.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.875.43'
$ws.Range('E3').Value = '  +2.36%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.48%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.49'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.005'
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4594'
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3860'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07858'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9960'
$ws.Range('E10').Value = '  +3.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.79'
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.891.13'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.991'
$ws.Range('E13').Value = '  +1.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.711'
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06952'
$ws.Range('E15').Value = '  +1.32%  '
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.006'
$ws.Range('E18').Value = '  +1.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.87'
$ws.Range('E19').Value = '  +1.13%  '
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.757.48'
$ws.Range('E21').Value = '  +2.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.279'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.03'
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.130'
$ws.Range('E24').Value = '  +2.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.121.13'
$ws.Range('E25').Value = '  +1.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.60'
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.24'
$ws.Range('E27').Value = '  +0.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.769'
$ws.Range('E28').Value = '  +0.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.967'
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.98'
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09326'
$ws.Range('E31').Value = '  +0.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9188'
$ws.Range('E32').Value = '  -2.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.305'
$ws.Range('E33').Value = '  +0.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.340'
$ws.Range('E34').Value = '  +1.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.328'
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05763'
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.154'
$ws.Range('E37').Value = '  +1.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02071'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.710'
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5641'
$ws.Range('E40').Value = '  +0.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1790'
$ws.Range('E41').Value = '  +1.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.904'
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07213'
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.79'
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5296'
$ws.Range('E45').Value = '  +0.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.143'
$ws.Range('E46').Value = '  +1.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.125'
$ws.Range('E47').Value = '  -1.06%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '113.61'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.825'
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.411'
$ws.Range('E50').Value = '  +3.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.005'
$ws.Range('E51').Value = '  +0.40%  '
